# Set the "Available" marker (X) for the "Shoulder" row (row 6, column C)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C6").Value = "X"

# Update selection to match the author's recorded cursor position
$ws.Range("C10").Select()
